$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "MSG: None`n`nMSG: No decision was made regarding the movie for Friday.`n"
$ws.Range("D2").Value = "no_decision, "
$ws.Range("C3").Value = "MSG: None`n`nMSG: The decision has been recorded: `"Oppenheimer`" has been selected as the movie to be shown on Friday.`n"
$ws.Range("C4").Value = "MSG: None`n`nMSG: The decision to acquire the rights for both movies has been successfully recorded.`n"
$ws.Range("D4").Value = "both_movies, "
$ws.Range("C5").Value = "MSG: None`n`nMSG: The rights for the movie `"Barbie`" have been acquired.`n"
$ws.Range("C6").Value = "MSG: None`n`nMSG: I have recorded the decision to acquire the rights for `"Barbie`" to be shown on Friday.`n"
$ws.Range("C7").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday has not been made.`n"
$ws.Range("D7").Value = "no_decision, "
$ws.Range("C8").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday was not made.`n"
$ws.Range("D8").Value = "no_decision, "
$ws.Range("C9").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie was selected.`n"
$ws.Range("D9").Value = "no_decision, "
$ws.Range("C10").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie being chosen for Friday.`n"
$ws.Range("D10").Value = "no_decision, "
$ws.Range("C11").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie`" to be shown on Friday.`n"
$ws.Range("C12").Value = "MSG: None`n`nMSG: The decision resulted in no clear agreement on which movie to show on Friday, so the outcome is recorded as no decision.`n"
$ws.Range("D12").Value = "no_decision, "
$ws.Range("C13").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision was made about which movie to play on Friday.`n"
$ws.Range("D13").Value = "no_decision, "
$ws.Range("C14").Value = "MSG: None`n`nMSG: I have recorded the decision: no movie will be shown on Friday.`n"
$ws.Range("D14").Value = "no_decision, "
$ws.Range("C15").Value = "MSG: None`n`nMSG: The decision to choose a movie for Friday has not been made.`n"
$ws.Range("D15").Value = "no_decision, "
$ws.Range("C16").Value = "MSG: None`n`nMSG: The rights to `"Barbie`" have been successfully acquired for the upcoming showing.`n"
$ws.Range("C17").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie selected for Friday.`n"
$ws.Range("D17").Value = "no_decision, "
$ws.Range("C18").Value = "MSG: None`n`nMSG: The decision to acquire both movies' rights has been successfully recorded.`n"
$ws.Range("D18").Value = "both_movies, "
$ws.Range("C19").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie being selected for Friday.`n"
$ws.Range("D19").Value = "no_decision, "
$ws.Range("C20").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision.`n"
$ws.Range("D20").Value = "no_decision, "
$ws.Range("C21").Value = "MSG: None`n`nMSG: The decision has been recorded as no agreement was made regarding the movie for Friday.`n"
$ws.Range("D21").Value = "no_decision, "
$ws.Range("C22").Value = "MSG: None`n`nMSG: The decision reflects that no movie was selected for Friday.`n"
$ws.Range("D22").Value = "no_decision, "
$ws.Range("C23").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday was not made.`n"
$ws.Range("D23").Value = "no_decision, "
$ws.Range("C24").Value = "MSG: None`n`nMSG: The decision for Friday's movie remains unresolved.`n"
$ws.Range("D24").Value = "no_decision, "
$ws.Range("C25").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision.`"`n"
$ws.Range("D25").Value = "no_decision, "
$ws.Range("C26").Value = "MSG: None`n`nMSG: I have recorded the decision to acquire the rights for `"Barbie.`"`n"
$ws.Range("C27").Value = "MSG: None`n`nMSG: The decision has been recorded: `"Barbie`" will be acquired for Friday's screening.`n"
$ws.Range("C28").Value = "MSG: None`n`nMSG: The decision has been recorded, and no movie was determined for Friday.`n"
$ws.Range("D28").Value = "no_decision, "
$ws.Range("C29").Value = "MSG: None`n`nMSG: The decision-making process concluded without a selection for Friday’s movie.`n"
$ws.Range("D29").Value = "no_decision, "
$ws.Range("C30").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie being selected.`n"
$ws.Range("D30").Value = "no_decision, "
$ws.Range("C31").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision about Friday's movie can be made.`n"
$ws.Range("D31").Value = "no_decision, "
$ws.Range("C32").Value = "MSG: None`n`nMSG: The committee did not reach a decision about which movie to play on Friday.`n"
$ws.Range("D32").Value = "no_decision, "
$ws.Range("C33").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie.`"`n"
$ws.Range("C34").Value = "MSG: None`n`nMSG: The decision-making process did not result in an agreement about which movie to show on Friday.`n"
$ws.Range("D34").Value = "no_decision, "
$ws.Range("C35").Value = "MSG: None`n`nMSG: The rights to both movies have been acquired for the screening on Friday.`n"
$ws.Range("C36").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights to `"Oppenheimer.`"`n"
$ws.Range("C37").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision was made about the movie to be shown on Friday.`n"
$ws.Range("D37").Value = "no_decision, "
$ws.Range("C38").Value = "MSG: None`n`nMSG: The function has been called, indicating that no decision about Friday's movie was made.`n"
$ws.Range("D38").Value = "no_decision, "
$ws.Range("C39").Value = "MSG: None`n`nMSG: The decision-making process ended without a selection for Friday's movie, resulting in no movie being chosen.`n"
$ws.Range("D39").Value = "no_decision, "
$ws.Range("C40").Value = "MSG: None`n`nMSG: The decision from the committee resulted in no movie being selected for Friday's showing.`n"
$ws.Range("D40").Value = "no_decision, "
$ws.Range("C41").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Barbie.`"`n"
$ws.Range("C42").Value = "MSG: None`n`nMSG: The decision-making process resulted in no consensus about which movie to show on Friday. Therefore, the outcome is recorded as `"no decision.`"`n"
$ws.Range("D42").Value = "no_decision, "
$ws.Range("C43").Value = "MSG: None`n`nMSG: The decision-making process did not lead to a clear choice of movie, resulting in no decision being made.`n"
$ws.Range("D43").Value = "no_decision, "
$ws.Range("C44").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie.`"`n"
$ws.Range("C45").Value = "MSG: None`n`nMSG: The decision was made to not select a movie for Friday, as there was no consensus reached between the committee members.`n"
$ws.Range("D45").Value = "no_decision, "
$ws.Range("C46").Value = "MSG: None`n`nMSG: The decision process resulted in no agreement on which movie to show on Friday, so no movie was selected.`n"
$ws.Range("D46").Value = "no_decision, "
$ws.Range("C47").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie selection was made for Friday.`n"
$ws.Range("D47").Value = "no_decision, "
$ws.Range("C48").Value = "MSG: None`n`nMSG: The committee did not reach a decision regarding which movie to show on Friday.`n"
$ws.Range("D48").Value = "no_decision, "
$ws.Range("C49").Value = "MSG: None`n`nMSG: The decision-making process concluded without a clear agreement on which movie to show on Friday, so no decision has been reached regarding the acquisition of any movie rights.`n"
$ws.Range("D49").Value = "no_decision, "
$ws.Range("C50").Value = "MSG: None`n`nMSG: The decision has been recorded, and no movie was selected for Friday's showing.`n"
$ws.Range("D50").Value = "no_decision, "
$ws.Range("C51").Value = "MSG: None`n`nMSG: The rights for both movies have been acquired.`n"
$ws.Range("C52").Value = "MSG: None`n`nMSG: The decision has been recorded and no movie will be shown on Friday.`n"
$ws.Range("D52").Value = "no_decision, "
$ws.Range("C53").Value = "MSG: None`n`nMSG: The decision regarding the movie to show on Friday has ended without reaching a consensus, so no movie rights will be acquired at this time.`n"
$ws.Range("D53").Value = "no_decision, "
$ws.Range("C54").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie being selected for Friday.`n"
$ws.Range("D54").Value = "no_decision, "
$ws.Range("C55").Value = "MSG: None`n`nMSG: The decision has been recorded as no selection was made regarding the movie for Friday.`n"
$ws.Range("D55").Value = "no_decision, "
$ws.Range("C56").Value = "MSG: None`n`nMSG: The decision has been recorded as a no decision regarding the movie to be shown on Friday.`n"
$ws.Range("D56").Value = "no_decision, "
$ws.Range("C57").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been made.`n"
$ws.Range("C58").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday has not been made.`n"
$ws.Range("D58").Value = "no_decision, "
$ws.Range("C59").Value = "MSG: None`n`nMSG: The decision regarding the movie to show on Friday has not been made, so there is no movie to acquire rights for at this time.`n"
$ws.Range("D59").Value = "no_decision, "
$ws.Range("C60").Value = "MSG: None`n`nMSG: The movie `"Barbie`" has been successfully selected for the event on Friday.`n"
$ws.Range("C61").Value = "MSG: None`n`nMSG: The conversation ended without a decision about which movie to show on Friday.`n"
$ws.Range("D61").Value = "no_decision, "
$ws.Range("C62").Value = "MSG: None`n`nMSG: The decision process has concluded with no movie selected for Friday.`n"
$ws.Range("D62").Value = "no_decision, "
$ws.Range("C63").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision`" regarding the movie to be shown on Friday.`n"
$ws.Range("D63").Value = "no_decision, "
$ws.Range("C64").Value = "MSG: None`n`nMSG: The conversation has not led to a decision about what movie will be shown on Friday, so I will record that as a 'no decision.'`n"
$ws.Range("D64").Value = "no_decision, "
$ws.Range("C65").Value = "MSG: None`n`nMSG: I have successfully acquired the rights for both movies, `"Oppenheimer`" and `"Barbie.`"`n"
$ws.Range("C66").Value = "MSG: None`n`nMSG: The function has been called, and the decision process confirmed that there was no agreement on a movie for Friday.`n"
$ws.Range("D66").Value = "no_decision, "
$ws.Range("C67").Value = "MSG: None`n`nMSG: The decision has been recorded that no movie was selected for Friday.`n"
$ws.Range("D67").Value = "no_decision, "
